$wb = $excel.ActiveWorkbook

# --- 1. Add new (incremental) hidden filter-database defined name on "template" sheet ---
$wsTemplate = $wb.Worksheets.Item("template")
$wsTemplate.Names.Add("_xlnm._FilterDatabase_0_0_0_0_0", "=template!`$A`$1:`$V`$64")

# --- 2. Fix typo in shared text used by rows 38-64, column A ---
$wsTemplate.Cells.Replace("Set FRIA Modreloramas", "Set FRIA Moderloramas")

# --- 3. Update sheet view: scroll position + active selection ---
$wsTemplate.Activate()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$wsTemplate.Range("A38").Select()

# --- 4. Column width tweaks on "template" sheet (slightly narrower columns) ---
$wsTemplate.Columns.Item(1).ColumnWidth = 52.666666666666664
$wsTemplate.Columns.Item(2).ColumnWidth = 34.5
$wsTemplate.Columns.Item(3).ColumnWidth = 31.166666666666668
$wsTemplate.Range($wsTemplate.Cells.Item(1,4), $wsTemplate.Cells.Item(1,15)).EntireColumn.ColumnWidth = 7.333333333333333
$wsTemplate.Columns.Item(16).ColumnWidth = 31.666666666666668
$wsTemplate.Columns.Item(17).ColumnWidth = 31.833333333333332
$wsTemplate.Columns.Item(18).ColumnWidth = 7.333333333333333
$wsTemplate.Columns.Item(19).ColumnWidth = 30.166666666666668
$wsTemplate.Range($wsTemplate.Cells.Item(1,20), $wsTemplate.Cells.Item(1,21)).EntireColumn.ColumnWidth = 7.333333333333333
$wsTemplate.Range($wsTemplate.Cells.Item(1,22), $wsTemplate.Cells.Item(1,1025)).EntireColumn.ColumnWidth = 6.666666666666667

# --- 5. Row heights 38-57 and 59-63 bump from 12.75 to 12.8 ---
for ($r = 38; $r -le 57; $r++) {
    $wsTemplate.Rows.Item($r).RowHeight = 12.8
}
for ($r = 59; $r -le 63; $r++) {
    $wsTemplate.Rows.Item($r).RowHeight = 12.8
}

# --- 6. Add explicit column formatting (single band) to the other two sheets ---
$wsHoja2 = $wb.Worksheets.Item("Hoja2")
$wsHoja2.Range($wsHoja2.Cells.Item(1,1), $wsHoja2.Cells.Item(1,1025)).EntireColumn.ColumnWidth = 7.666666666666667

$wsStoreTypes = $wb.Worksheets.Item("store types")
$wsStoreTypes.Range($wsStoreTypes.Cells.Item(1,1), $wsStoreTypes.Cells.Item(1,1025)).EntireColumn.ColumnWidth = 7.666666666666667
